$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Ativação:" date text (row 8) ---------------------------------------
# NOTE: the original cell stores the literal text "01/01/2020" (a shared
# string), not a real date. Assigning a date-looking literal straight to
# .Value would make Excel auto-convert it to a date serial number (and
# fork a new number-format style), which is NOT what the source workbook
# does. To keep it as plain text we build the string with a formula in a
# scratch cell, copy it, and paste *values only* into the target cells -
# PasteSpecial(xlPasteValues) carries over the computed text without
# triggering Excel's "looks like a date" literal-entry heuristic, and
# without touching the destination cell's existing style.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = '="01/01/" & "2022"'
$scratch.Copy()
$ws.Range("B8").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163)   # xlPasteValues
$scratch.ClearContents()
$excel.CutCopyMode = 0

# --- "Programa resumido:" / "Short syllabus:" (row 15/16) ----------------
$ws.Range("B15").Value = "Geomorfologia Fluvial; Padrões de Drenagem; Escoamentos hidráulicos; medidores; bocais; instrumentos de medição"
$ws.Range("C15").Value = "Geomorfologia Fluvial; Padrões de Drenagem; Escoamentos hidráulicos; medidores; bocais; instrumentos de medição"

$ws.Range("B16").Value = "River Geomorphology; Drainage Patterns; Hydraulic flow; meters; nozzles; measuring instruments."
$ws.Range("C16").Value = "River Geomorphology; Drainage Patterns; Hydraulic flow; meters; nozzles; measuring instruments."

# --- "Programa:" / "Syllabus:" (row 17/18) --------------------------------
$ws.Range("B17").Value = "- As teorias geomorfológicas;- Processos e Formas do relevo;- Processos fluviais, morfologias fluviais e padrões de drenagem;- Precipitação;- Infiltração;- Evapotranspiração;- Escoamento superficial;- Instrumentos de medição (Calhas, vertedores e registros);- Operação de reservatórios;- Vazões máximas e mínimas: distribuição de frequência, hidrograma unitário.- Água subterrânea, aquíferos e poços;"
$ws.Range("C17").Value = "- As teorias geomorfológicas;- Processos e Formas do relevo;- Processos fluviais, morfologias fluviais e padrões de drenagem;- Precipitação;- Infiltração;- Evapotranspiração;- Escoamento superficial;- Instrumentos de medição (Calhas, vertedores e registros);- Operação de reservatórios;- Vazões máximas e mínimas: distribuição de frequência, hidrograma unitário.- Água subterrânea, aquíferos e poços;"

$ws.Range("B18").Value = "- Geomorphological theories;- Processes and Forms of relief;- River processes, river morphologies and drainage patterns;- Precipitation;- Infiltration;- Evapotranspiration;- Surface runoff;- Measuring instruments (gutters, spillways and registers);- Reservoir operation;- Maximum and minimum flow rates: frequency distribution, unit hydrograph.- Groundwater, aquifers and wells;"
$ws.Range("C18").Value = "- Geomorphological theories;- Processes and Forms of relief;- River processes, river morphologies and drainage patterns;- Precipitation;- Infiltration;- Evapotranspiration;- Surface runoff;- Measuring instruments (gutters, spillways and registers);- Reservoir operation;- Maximum and minimum flow rates: frequency distribution, unit hydrograph.- Groundwater, aquifers and wells;"
